# Append two new daily rows (row 119 = 2025-11-27 serial 45988,
# row 120 = 2025-11-28 serial 45989) to each of the 6 worksheets,
# extending the used range from A1:B118 to A1:B120.
#
# Column A keeps the same date number-format as the existing rows
# (copied from the last existing row so the style index is reused
# rather than a new one being created), column B is a plain number.

$wb = $excel.ActiveWorkbook

# sheetName -> [row119 remn_amt, row120 remn_amt]
$newData = @{
    "LG생활건강"   = @(509882, 0)
    "아모레퍼시픽" = @(320040, 0)
    "한국콜마"     = @(128780, 0)
    "코스맥스"     = @(223176, 0)
    "에이피알"     = @(761281, 0)
    "달바글로벌"   = @(68956, 0)
}

$newDates = @(45988, 45989)

foreach ($ws in $wb.Worksheets) {
    $values = $newData[$ws.Name]
    if ($values -eq $null) { continue }

    $dateFormat = $ws.Cells.Item(118, 1).NumberFormat

    for ($i = 0; $i -lt 2; $i++) {
        $row = 119 + $i

        $dateCell = $ws.Cells.Item($row, 1)
        $dateCell.Value = $newDates[$i]
        $dateCell.NumberFormat = $dateFormat

        $amtCell = $ws.Cells.Item($row, 2)
        $amtCell.Value = $values[$i]
    }
}
